$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.355.68"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").Value = "'3.418.54"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'581.97"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("D6").Value = "'176.29"
$ws.Range("E6").Value = "  -2.27%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'3.411.45"
$ws.Range("E8").Value = "  +1.01%  "

# Row 9
$ws.Range("E9").Value = "  -0.71%  "

# Row 10
$ws.Range("E10").Value = "  +0.75%  "

# Row 11
$ws.Range("E11").Value = "  -1.33%  "

# Row 12
$ws.Range("D12").Value = "'48.63"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("E13").Value = "  -1.97%  "

# Row 14
$ws.Range("D14").Value = "'690.38"
$ws.Range("E14").Value = "  +0.68%  "

# Row 15
$ws.Range("D15").Value = "'3.967.24"
$ws.Range("E15").Value = "  +1.00%  "

# Row 17
$ws.Range("D17").Value = "'69.435.32"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18
$ws.Range("D18").Value = "'3.421.18"
$ws.Range("E18").Value = "  +1.14%  "

# Row 19
$ws.Range("E19").Value = "  +0.85%  "

# Row 20
$ws.Range("D20").Value = "'17.59"
$ws.Range("E20").Value = "  -0.72%  "

# Row 21
$ws.Range("D21").Value = "'11.32"
$ws.Range("E21").Value = "  -0.34%  "

# Row 22
$ws.Range("D22").Value = "'0.893"
$ws.Range("E22").Value = "  -0.83%  "

# Row 23
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("D24").Value = "'16.84"
$ws.Range("E24").Value = "  -1.27%  "

# Row 25
$ws.Range("D25").Value = "'100.96"
$ws.Range("E25").Value = "  -3.50%  "

# Row 26
$ws.Range("E26").Value = "  -1.07%  "

# Row 27
$ws.Range("E27").Value = "  -2.14%  "

# Row 28
$ws.Range("D28").Value = "'9.52"
$ws.Range("E28").Value = "  -0.89%  "

# Row 29
$ws.Range("D29").Value = "'33.36"
$ws.Range("E29").Value = "  -2.82%  "

# Row 30
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").Value = "'7.02"
$ws.Range("E31").Value = "  +0.86%  "

# Row 32
$ws.Range("D32").Value = "'576.31"
$ws.Range("E32").Value = "  +3.49%  "

# Row 33
$ws.Range("D33").Value = "'3.65"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("D34").Value = "'10.98"
$ws.Range("E34").Value = "  -1.83%  "

# Row 35
$ws.Range("D35").Value = "'58.29"
$ws.Range("E35").Value = "  +0.62%  "

# Row 36
$ws.Range("E36").Value = "  -3.25%  "

# Row 37
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("D38").Value = "'3.569.06"
$ws.Range("E38").Value = "  -3.57%  "

# Row 39
$ws.Range("D39").Value = "'0.139"
$ws.Range("E39").Value = "  -0.98%  "

# Row 40
$ws.Range("D40").Value = "'34.77"
$ws.Range("E40").Value = "  -0.23%  "

# Row 41
$ws.Range("D41").Value = "'0.0₃0726"
$ws.Range("E41").Value = "  +3.10%  "

# Row 42
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("E43").Value = "  -1.10%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.331"
$ws.Range("E44").Value = "  -2.27%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0416"
$ws.Range("E45").Value = "  -0.53%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'1.45"
$ws.Range("E46").Value = "  +4.46%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.64"
$ws.Range("E47").Value = "  -0.38%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.128"
$ws.Range("E48").Value = "  -1.26%  "

# Row 49
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.23%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'132.38"
$ws.Range("E50").Value = "  -0.19%  "

# Row 51
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "'2.62"
$ws.Range("E51").Value = "  +1.81%  "

# Clear the "quote prefix" style flag introduced by the leading apostrophe above
# so formatting stays identical to the original (unstyled) cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
